# "Added Priority to Login&Queue"
#
# The underlying edit is a text correction on the "Login" sheet: the
# message shown for missing-field validation ("Please fill out this
# field") gets a trailing period added, becoming "Please fill out this
# field." (which matches the wording already used elsewhere in the
# workbook). Updating every cell that held the old text lets the
# now-unused shared string drop out of the table on save.
#
# Along with that, the active/selected sheet moves from "Queue" to
# "Login", with the selection parked on D4.

$wb = $excel.ActiveWorkbook

$login = $wb.Sheets.Item("Login")
$login.Range("D2").Value = "Please fill out this field."
$login.Range("D3").Value = "Please fill out this field."
$login.Range("D4").Value = "Please fill out this field."

$login.Activate()
$login.Range("D4").Select()
